# Update gh-pages output data (generated at 456a3b4)
# Applies updated "想去人数" (interest count, column F) and "最低票价"
# (minimum ticket price, column G) figures to the three sheets that carry
# per-event rows: 展览 (Exhibitions), 演出 (Performances), and the combined
# 全部类型 (All types) rollup sheet. 本地生活 (Local life) is unaffected.

function Set-CellIfMatches {
    param($ws, $cellRef, $oldVal, $newVal)
    $cur = $ws.Range($cellRef).Value2
    if ($cur -ne $oldVal) {
        Write-Host ("WARN: " + $ws.Name + "!" + $cellRef + " expected " + $oldVal + " but found " + $cur)
    }
    $ws.Range($cellRef).Value2 = $newVal
}

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
Set-CellIfMatches $ws1 "F3" 3838 3841
Set-CellIfMatches $ws1 "G3" 70 80
Set-CellIfMatches $ws1 "G5" 75 80
Set-CellIfMatches $ws1 "F6" 3834 3835
Set-CellIfMatches $ws1 "G6" 30 75
Set-CellIfMatches $ws1 "F8" 201 203
Set-CellIfMatches $ws1 "F10" 8667 8682
Set-CellIfMatches $ws1 "F11" 489 490
Set-CellIfMatches $ws1 "F12" 81 82
Set-CellIfMatches $ws1 "F17" 94 96
Set-CellIfMatches $ws1 "F18" 362 363
Set-CellIfMatches $ws1 "F19" 10987 10994
Set-CellIfMatches $ws1 "F28" 2680 2683
Set-CellIfMatches $ws1 "F29" 2081 2082
Set-CellIfMatches $ws1 "F30" 42 43
Set-CellIfMatches $ws1 "F33" 899 901
Set-CellIfMatches $ws1 "F34" 4087 4090
Set-CellIfMatches $ws1 "F35" 2567 2568
Set-CellIfMatches $ws1 "F37" 2585 2587
Set-CellIfMatches $ws1 "F38" 3028 3029
Set-CellIfMatches $ws1 "F39" 1248 1249
Set-CellIfMatches $ws1 "F40" 174 175
Set-CellIfMatches $ws1 "F41" 755 756
Set-CellIfMatches $ws1 "F43" 331 334
Set-CellIfMatches $ws1 "F44" 46 47
Set-CellIfMatches $ws1 "F45" 111 112
Set-CellIfMatches $ws1 "F46" 128 129
Set-CellIfMatches $ws1 "F49" 85 86

$ws2 = $wb.Worksheets.Item("演出")
Set-CellIfMatches $ws2 "F19" 175 176
Set-CellIfMatches $ws2 "F22" 54 59

$ws4 = $wb.Worksheets.Item("全部类型")
Set-CellIfMatches $ws4 "F3" 3838 3841
Set-CellIfMatches $ws4 "G3" 70 80
Set-CellIfMatches $ws4 "G6" 75 80
Set-CellIfMatches $ws4 "F7" 3834 3835
Set-CellIfMatches $ws4 "G7" 30 75
Set-CellIfMatches $ws4 "F10" 201 203
Set-CellIfMatches $ws4 "F11" 8667 8682
Set-CellIfMatches $ws4 "F12" 489 490
Set-CellIfMatches $ws4 "F16" 94 96
Set-CellIfMatches $ws4 "F17" 362 363
Set-CellIfMatches $ws4 "F18" 10987 10994
Set-CellIfMatches $ws4 "F29" 2680 2683
Set-CellIfMatches $ws4 "F30" 2081 2082
Set-CellIfMatches $ws4 "F32" 899 901
Set-CellIfMatches $ws4 "F34" 4087 4090
Set-CellIfMatches $ws4 "F35" 2567 2568
Set-CellIfMatches $ws4 "F37" 2585 2587
Set-CellIfMatches $ws4 "F38" 3028 3029
Set-CellIfMatches $ws4 "F39" 54 59
Set-CellIfMatches $ws4 "F40" 1248 1249
Set-CellIfMatches $ws4 "F41" 174 175
Set-CellIfMatches $ws4 "F42" 755 756
Set-CellIfMatches $ws4 "F44" 331 334
Set-CellIfMatches $ws4 "F45" 111 112
Set-CellIfMatches $ws4 "F46" 128 129
Set-CellIfMatches $ws4 "F49" 85 86

Write-Host "Done applying gh-pages data refresh."
